$rowEdits = @(
    @{ Sheet="ALC"; Row=15; Set=@{"H"=1265.8572; "I"=1265.8572; "K"=3797.5716; "M"=-3628.5716}; Clear=@() },
    @{ Sheet="ALC"; Row=17; Set=@{"H"=986.7941; "J"=983.36206; "L"=2950.08618; "N"=-3286.08618}; Clear=@() },
    @{ Sheet="ALC"; Row=58; Set=@{"H"=1178.4286; "I"=791.5; "J"=3500; "K"=2374.5; "L"=10500; "M"=-2224.5; "N"=-10800}; Clear=@() },
    @{ Sheet="ALC"; Row=86; Set=@{"H"=15393.637; "I"=13258.889; "K"=13258.889; "M"=-12135.889}; Clear=@() },
    @{ Sheet="ALC"; Row=89; Set=@{"H"=15393.637; "I"=13258.889; "K"=66294.44499999999; "M"=-60678.44499999999}; Clear=@() },
    @{ Sheet="ALC"; Row=106; Set=@{"H"=8346.182000000001; "I"=8346.182000000001; "K"=8346.182000000001; "M"=-7715.182000000001}; Clear=@() },
    @{ Sheet="ALC"; Row=137; Set=@{"H"=3737.0164; "I"=1918.3208; "K"=5754.9624; "M"=-3204.9624}; Clear=@() },
    @{ Sheet="ALC"; Row=138; Set=@{"H"=2367.0444; "I"=1827.6052; "J"=5295.4287; "K"=5482.8156; "L"=15886.2861; "M"=-342.8155999999999; "N"=-26166.2861}; Clear=@() },
    @{ Sheet="ARM"; Row=74; Set=@{"H"=6355.8433; "I"=3490.9; "J"=16773.818; "K"=3490.9; "L"=16773.818; "M"=-2616.9; "N"=-18521.818}; Clear=@() },
    @{ Sheet="ARM"; Row=77; Set=@{"H"=6355.8433; "I"=3490.9; "J"=16773.818; "K"=17454.5; "L"=83869.09; "M"=-13086.5; "N"=-92605.09}; Clear=@() },
    @{ Sheet="CRP"; Row=23; Set=@{"H"=30000; "I"=0; "J"=30000; "K"=0; "L"=30000; "N"=-30480}; Clear=@("M") },
    @{ Sheet="CRP"; Row=27; Set=@{"H"=30000; "I"=0; "J"=30000; "K"=0; "L"=30000; "N"=-30384}; Clear=@("M") },
    @{ Sheet="CRP"; Row=31; Set=@{"H"=4783; "I"=3576.625; "K"=3576.625; "M"=-3281.625}; Clear=@() },
    @{ Sheet="CRP"; Row=34; Set=@{"H"=4783; "I"=3576.625; "K"=3576.625; "M"=-3374.625}; Clear=@() },
    @{ Sheet="CRP"; Row=122; Set=@{"H"=7706.8535; "I"=1629.9117; "K"=4889.7351; "M"=-2439.7351}; Clear=@() },
    @{ Sheet="CRP"; Row=134; Set=@{"H"=2003.9615; "I"=2004.16; "K"=6012.48; "M"=-3477.48}; Clear=@() },
    @{ Sheet="CUL"; Row=2; Set=@{"H"=34.384617; "I"=34.384617; "J"=0; "K"=206.307702; "L"=0; "M"=-93.30770200000001}; Clear=@("N") },
    @{ Sheet="CUL"; Row=9; Set=@{"H"=758024.5600000001; "I"=333367.34; "J"=788357.2; "K"=1000102.02; "L"=2365071.6; "M"=-999878.02; "N"=-2365519.6}; Clear=@() },
    @{ Sheet="CUL"; Row=10; Set=@{"H"=2874.75; "I"=750; "J"=4999.5; "K"=2250; "L"=14998.5; "M"=-2111; "N"=-15276.5}; Clear=@() },
    @{ Sheet="CUL"; Row=13; Set=@{"H"=161; "I"=0; "J"=161; "K"=0; "L"=483; "N"=-819}; Clear=@("M") },
    @{ Sheet="CUL"; Row=16; Set=@{"H"=0; "J"=0; "L"=0}; Clear=@("N") },
    @{ Sheet="CUL"; Row=22; Set=@{"H"=800.6667; "J"=2002; "L"=6006; "N"=-6344}; Clear=@() },
    @{ Sheet="CUL"; Row=27; Set=@{"H"=800.6667; "J"=2002; "L"=6006; "N"=-6210}; Clear=@() },
    @{ Sheet="CUL"; Row=33; Set=@{"H"=3840}; Clear=@() },
    @{ Sheet="CUL"; Row=50; Set=@{"H"=308.40475; "I"=301.17648; "J"=313.32; "K"=903.52944; "L"=939.96; "M"=-422.52944; "N"=-1901.96}; Clear=@() },
    @{ Sheet="CUL"; Row=53; Set=@{"H"=308.40475; "I"=301.17648; "J"=313.32; "K"=903.52944; "L"=939.96; "M"=-422.52944; "N"=-1901.96}; Clear=@() },
    @{ Sheet="CUL"; Row=81; Set=@{"H"=2104.1428; "J"=3375; "L"=10125; "N"=-12371}; Clear=@() },
    @{ Sheet="CUL"; Row=84; Set=@{"H"=2104.1428; "J"=3375; "L"=30375; "N"=-41607}; Clear=@() },
    @{ Sheet="CUL"; Row=94; Set=@{"H"=10987.25; "I"=4299.3335; "K"=12898.0005; "M"=-12222.0005}; Clear=@() },
    @{ Sheet="CUL"; Row=122; Set=@{"H"=1242333.8; "J"=1896.1818; "L"=17065.6362; "N"=-21965.6362}; Clear=@() },
    @{ Sheet="GSM"; Row=75; Set=@{"H"=49933.332; "J"=49933.332; "L"=49933.332; "N"=-51681.332}; Clear=@() },
    @{ Sheet="GSM"; Row=78; Set=@{"H"=49933.332; "J"=49933.332; "L"=149799.996; "N"=-158535.996}; Clear=@() },
    @{ Sheet="LTW"; Row=22; Set=@{"H"=4124.125; "J"=4142; "L"=4142; "N"=-4732}; Clear=@() },
    @{ Sheet="LTW"; Row=27; Set=@{"H"=4124.125; "J"=4142; "L"=4142; "N"=-4356}; Clear=@() },
    @{ Sheet="LTW"; Row=34; Set=@{"H"=0; "I"=0; "K"=0}; Clear=@("M") },
    @{ Sheet="LTW"; Row=61; Set=@{"H"=7112.9062; "I"=8133.1665; "J"=4052.125; "K"=8133.1665; "L"=4052.125; "M"=-7931.1665; "N"=-4456.125}; Clear=@() },
    @{ Sheet="LTW"; Row=68; Set=@{"H"=12275.167; "I"=10828.833; "K"=10828.833; "M"=-10079.833}; Clear=@() },
    @{ Sheet="LTW"; Row=71; Set=@{"H"=12275.167; "I"=10828.833; "K"=54144.165; "M"=-50400.165}; Clear=@() },
    @{ Sheet="LTW"; Row=99; Set=@{"H"=42347.25; "I"=43296.332; "J"=39500; "K"=43296.332; "L"=39500; "M"=-40301.332; "N"=-45490}; Clear=@() },
    @{ Sheet="LTW"; Row=113; Set=@{"H"=7112.9062; "I"=8133.1665; "J"=4052.125; "K"=8133.1665; "L"=4052.125; "M"=-5963.1665; "N"=-8392.125}; Clear=@() },
    @{ Sheet="LTW"; Row=132; Set=@{"H"=2982.348; "I"=1978.3846; "K"=5935.1538; "M"=-3405.1538}; Clear=@() },
    @{ Sheet="LTW"; Row=136; Set=@{"H"=1789.5483; "I"=1505.25; "J"=4443; "K"=4515.75; "L"=13329; "M"=-1965.75; "N"=-18429}; Clear=@() },
    @{ Sheet="WVR"; Row=81; Set=@{"H"=171491.67; "I"=7500; "K"=15000; "M"=-13939}; Clear=@() },
    @{ Sheet="WVR"; Row=84; Set=@{"H"=171491.67; "I"=7500; "K"=75000; "M"=-69696}; Clear=@() },
    @{ Sheet="WVR"; Row=113; Set=@{"H"=2570.75; "I"=2570.75; "J"=0; "K"=7712.25; "L"=0; "M"=-5542.25}; Clear=@("N") }
)

$wb = $excel.ActiveWorkbook

foreach ($edit in $rowEdits) {
    $ws = $wb.Worksheets.Item($edit.Sheet)
    foreach ($col in $edit.Set.Keys) {
        $addr = "$col$($edit.Row)"
        $ws.Range($addr).Value = $edit.Set[$col]
    }
    foreach ($col in $edit.Clear) {
        $addr = "$col$($edit.Row)"
        $ws.Range($addr).Value = $null
    }
}
